$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value in A3 (next row after existing data)
$ws.Range("A3").Value = 3

# Move/update the selection to N5
$ws.Range("N5").Select()
